{"js": "// Table 21 (\"List of Institutional factors.\") stores the ORGANIZATION\n// column as a raw ~8-decimal fraction (e.g. \"0.03225806\") or the literal\n// \"NA\". Round the displayed figure to 2 decimals (and normalise the\n// missing-value marker to \"-\") without touching any other column.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Locate the right table by its header row rather than assuming index 0.\nlet target = null;\nlet orgColIndex = -1;\nfor (const table of tables.items) {\n  table.load(\"values\");\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  const header = table.values && table.values.length > 0 ? table.values[0] : [];\n  const idx = header.indexOf(\"ORGANIZATION\");\n  if (idx !== -1) {\n    target = table;\n    orgColIndex = idx;\n    break;\n  }\n}\n\nif (target) {\n  const rowCount = target.values.length;\n\n  // Read every cell in the ORGANIZATION column up front.\n  const cells = [];\n  for (let r = 1; r < rowCount; r++) {\n    const cell = target.getCell(r, orgColIndex);\n    cell.load(\"value\");\n    cells.push(cell);\n  }\n  await context.sync();\n\n  // Work out (and apply) the new text for each cell.\n  for (const cell of cells) {\n    const raw = (cell.value || \"\").trim();\n    let next = null;\n\n    if (raw === \"NA\") {\n      next = \"-\";\n    } else if (raw !== \"\" && !isNaN(Number(raw))) {\n      const rounded = Number(raw).toFixed(2);\n      if (rounded !== raw) {\n        next = rounded;\n      }\n    }\n\n    if (next !== null) {\n      cell.body.getRange().insertText(next, Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Table 21 (\"List of Institutional factors.\") stores the ORGANIZATION\n# column as a raw ~8-decimal fraction (e.g. \"0.03225806\") or the literal\n# \"NA\". Round the displayed figure to 2 decimals (and normalise the\n# missing-value marker to \"-\") without touching any other column.\n\n$d = $word.ActiveDocument\n\nfunction Clean-CellText($text) {\n    return $text -replace \"[\\x07\\x0d]\", \"\"\n}\n\nfor ($t = 1; $t -le $d.Tables.Count; $t++) {\n    $tbl = $d.Tables.Item($t)\n\n    # Find the ORGANIZATION column from the header row (row 1).\n    $orgCol = -1\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        $headerText = Clean-CellText $tbl.Cell(1, $c).Range.Text\n        if ($headerText.Equals(\"ORGANIZATION\")) {\n            $orgCol = $c\n            break\n        }\n    }\n\n    if ($orgCol -eq -1) {\n        continue\n    }\n\n    for ($r = 2; $r -le $tbl.Rows.Count; $r++) {\n        $cell = $tbl.Cell($r, $orgCol)\n        $raw = (Clean-CellText $cell.Range.Text).Trim()\n\n        $next = $null\n        if ($raw.Equals(\"NA\")) {\n            $next = \"-\"\n        } elseif ($raw -match '^[0-9]+(\\.[0-9]+)?$') {\n            $num = [double]$raw\n            $rounded = $num.ToString(\"0.00\")\n            if (-not $rounded.Equals($raw)) {\n                $next = $rounded\n            }\n        }\n\n        if ($null -ne $next) {\n            $cell.Range.Text = $next\n        }\n    }\n}\n"}
